$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used with PasteSpecial below.
$xlPasteFormats = -4122

# --- Column P, row 3 (thin bottom-border-only row, blank) ---
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial($xlPasteFormats)

# --- Column P, row 4 (year header) ---
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial($xlPasteFormats)
$ws.Range("P4").Value2 = 2021

# --- Column P, row 5 (totals row) - style 18 (matches O8's style) ---
$ws.Range("O8").Copy()
$ws.Range("P5").PasteSpecial($xlPasteFormats)
$ws.Range("P5").Value2 = 9038

# --- Column P, row 6 (blank sub-header row) ---
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial($xlPasteFormats)

# --- Column P, row 7 - style 17 (matches O6's blank style) ---
$ws.Range("O6").Copy()
$ws.Range("P7").PasteSpecial($xlPasteFormats)
$ws.Range("P7").Value2 = 8587

# --- Column P, row 8 - style 17 (matches O6's blank style) ---
$ws.Range("O6").Copy()
$ws.Range("P8").PasteSpecial($xlPasteFormats)
$ws.Range("P8").Value2 = 451

# --- Column P, row 9 (blank sub-header row) ---
$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial($xlPasteFormats)

# --- Column P, rows 10-25: "..." (no data) marker, same style as column O ---
for ($r = 10; $r -le 25; $r++) {
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial($xlPasteFormats)
    $ws.Range("P$r").Value2 = $ws.Range("O$r").Value2
}

# Clear clipboard marching ants / copy mode.
$excel.CutCopyMode = $false

# Restore the active selection to Q4, matching the saved view state.
$ws.Range("Q4").Select() | Out-Null
